# The underlying data (daily price observations for Kiwi / Agrícola del
# Norte S.A. de Arica) was refreshed: rows 2-30 keep their identifying
# columns (Mercado/Región/Codreg/Tipo/Producto/Categoría) but the
# observation-specific columns (Fecha, Variedad, Calidad, Volumen, Precios,
# Unidad, Origen, Precio $/Kg, Kg/unidad) are reassigned across rows
# (a weekly re-pull reshuffled which date lines up with which row).
#
# Columns that move: D, K, L, M, N, O, P, Q, R, S, T
# Row mapping (new row <- old row, i.e. row R's new content equals the
# snapshot of row Map[R] taken before any writes happen):

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
  2  = 26
  3  = 20
  4  = 21
  5  = 25
  6  = 5
  7  = 28
  8  = 16
  9  = 30
  10 = 15
  11 = 7
  12 = 4
  13 = 17
  14 = 23
  15 = 22
  16 = 8
  17 = 3
  18 = 24
  19 = 12
  20 = 6
  21 = 19
  22 = 14
  23 = 9
  24 = 2
  25 = 11
  26 = 18
  27 = 29
  28 = 13
  29 = 27
  30 = 10
}

$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # D,K,L,M,N,O,P,Q,R,S,T

# Snapshot every source row's values for the moving columns before any
# writes occur, since several rows source from each other (cycles).
$snapshot = @{}
foreach ($r in 2..30) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in 2..30) {
    $src = $map[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c]
    }
}
